$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 14:22"

# Row 8 - Alemania
$ws.Range("B8").Value = 130214
$ws.Range("C8").Value = 142
$ws.Range("E8").Value = 58811
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 3203

# Row 14 - Paises Bajos
$ws.Range("B14").Value = 27419
$ws.Range("C14").Value = 868
$ws.Range("E14").Value = 24224
$ws.Range("G14").Value = 122
$ws.Range("H14").Value = 2945

# Row 20 - Austria
$ws.Range("B20").Value = 14146
$ws.Range("C20").Value = 105
$ws.Range("E20").Value = 6129

# Row 22 - Suecia
$ws.Range("B22").Value = 11445
$ws.Range("C22").Value = 497
$ws.Range("E22").Value = 10031
$ws.Range("F22").Value = 915
$ws.Range("G22").Value = 114
$ws.Range("H22").Value = 1033

# Row 62 - Croacia
$ws.Range("B62").Value = 1704
$ws.Range("C62").Value = 54
$ws.Range("D62").Value = 415
$ws.Range("E62").Value = 1258
$ws.Range("G62").Value = 6
$ws.Range("H62").Value = 31

# Row 78 - Republica de Macedonia
$ws.Range("D78").Value = 86
$ws.Range("E78").Value = 778
$ws.Range("G78").Value = 6
$ws.Range("H78").Value = 44

# Row 109 - Estado de Palestina
$ws.Range("D109").Value = 62
$ws.Range("E109").Value = 244
